# Update "Generate Report for Handback" timestamps.
# These cells hold plain text timestamps (not Excel date values),
# stored in the shared-strings table, so assigning a matching
# "yyyy-mm-dd HH:mm:ss" text string preserves their text type.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-30 21:14:38"

# zh-cn sheet: Correspond Handoff / Handback DateTime for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-30 21:14:33"
$wsZhCn.Range("K2").Value = "2016-08-30 21:14:50"

# de-de sheet: Correspond Handoff / Handback DateTime for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-30 21:14:38"
$wsDeDe.Range("K2").Value = "2016-08-30 21:14:57"
